# Edit script for UC02 - Manter Cliente.xlsx
# Applies text corrections to the 4 "Fluxo" use-case sheets, adjusts a couple
# of row heights, a vertical-alignment style tweak, and sheet1 page margins,
# matching the authors' revision of the use case documentation.

$wb = $excel.ActiveWorkbook

$sheets = @{}
$sheets[1] = $wb.Worksheets.Item(1)
$sheets[2] = $wb.Worksheets.Item(2)
$sheets[3] = $wb.Worksheets.Item(3)
$sheets[4] = $wb.Worksheets.Item(4)

# ---- Text updates (shared string contents) ----
$text17 = @"
P1 - O ator seleciona a opção "Clientes";
P2 - O sistema exibe a página de clientes;
P3 - O ator seleciona a opção "Cadastrar Cliente";
P4 - O sistema exibe a página de cadastro de cliente;
P5 - O ator preenche o formulário de cadastro da página 
com os dados necessários;
P6 - O ator seleciona a opção "Efetuar Cadastro";
P7 - O sistema exibe uma janela com as informações do cadastro 
para o ator confirmar;
P8 - O ator clica na opção "Confimar"
P9 - O sistema grava os dados do formulário no banco de dados;
P10 - O sistema exibe uma janela confirmando o sucesso do cadastro;
P11 - O ator clica na opção "OK";
P12 - O sistema exibe a página de cadastro de cliente;
P13 - Caso de Uso concluído.
"@
$sheets[1].Range("C10").Value = $text17

$text29 = @"
P1 - O ator seleciona a opção "Clientes";
P2 - O sistema exibe a página de clientes. Nesta página, há um botão de 
"Cadastrar Cliente" e uma lista de clientes com cada registro tendo 
os botões de: detalhes, alterar e excluir; 
P3 - Caso de Uso concluído.
"@
$sheets[2].Range("C10").Value = $text29

$text30 = @"
P2.1 - O ator pode escolher diferentes filtros para listar os
clientes;
P2.2 -  O ator pode pesquisar clientes específicos.
"@
$sheets[2].Range("C11").Value = $text30

$text31 = @"
P2.1 - Se o banco de dados não estiver acessível, o sistema exibe a 
mensagem de erro "Falha na comunicação com o banco de dados";
P2.2 - Se ocorrer erro ao executar o SQL, o sistema exibe a mensagem
de erro "Falha ao gravar os dados no banco de dados".
"@
$sheets[2].Range("C12").Value = $text31

$text37 = @"
P1 - O ator seleciona a opção "Clientes";
P2 - O sistema exibe a página de clientes;
P3 - O ator seleciona o cliente que quer alterar;
P4 - O sistema exibe uma página para exibir as informações e permitir 
alterações daquele cliente em específico;
P5 - O ator altera as informações do formulário;
P6 - O ator clica na opção "Alterar dados";
P7 - O sistema exibe uma janela com as informações alteradas do cliente 
para o ator confirmar;
P8 - O ator clica na opção "Confimar";
P9 - O sistema grava os dados do formulário no banco de dados;
P10 - O sistema exibe uma janela confirmando o sucesso da alteração;
P11 - O ator clica na opção "OK";
P12 - O sistema exibe a página de clientes;
P13 - Caso de Uso concluído.
"@
$sheets[3].Range("C10").Value = $text37

$text38 = @"
P6.1 - O ator clica em "Voltar", cancelando a operação 
de alteração;
P7.1 - O sistema avisa o ator caso algum campo obrigatório
não esteja preenchido;
P8.1 - O ator clica em "Cancelar", saindo da janela para voltar
o preenchimento do formulário.
"@
$sheets[3].Range("C11").Value = $text38

$text39 = @"
P4.1 - Se o banco de dados não estiver acessível, o sistema exibe a 
mensagem de erro "Falha na comunicação com o banco de dados";
P7.1 - Caso algum campo seja preenchido com valores inválidos, o 
sistema emite a mensagem "Formato de dados inválido. Preencha os 
campos corretamente."
P9.1 - Se ocorrer erro ao executar o SQL, o sistema exibe a mensagem
de erro "Falha ao gravar os dados no banco de dados"
P10.1 -  Se já existir algum cliente com algum dos campos únicos já 
cadastrado, o sistema exibe uma mensagem de erro "Algum campo único 
já está cadastrado no sistema. Tente novamente."
"@
$sheets[3].Range("C12").Value = $text39

$text44 = @"
P1 - O ator seleciona a opção "Clientes";
P2 - O sistema exibe a página de clientes;
P3 - O ator clica no botão "Excluir" de um cliente específico;
P4 - O sistema emite uma mensagem pedindo confirmação dessa ação;
P5 - O ator clica no botão "Confirmar";
P6 - O sistema atualiza o banco de dados;
P7 - Caso de Uso concluído.
"@
$sheets[4].Range("C10").Value = $text44

$text45 = @"
P5.1 - O ator clica no botão "Cancelar", cancelando a operação
de exclusão.
"@
$sheets[4].Range("C11").Value = $text45

$text46 = @"
P6.1 - Se o cliente tiver transações associadas que impeçam a exclusão, o sistema exibe a mensagem "O cliente possui transações ativas e não pode ser excluído.";
P6.2 - Se o banco de dados não estiver acessível, o sistema exibe a 
mensagem de erro "Falha na comunicação com o banco de dados";
P6.3 - Se ocorrer erro ao executar o SQL, o sistema exibe a mensagem
de erro "Falha ao excluir o cliente no banco de dados."
"@
$sheets[4].Range("C12").Value = $text46

# ---- Row height fix-ups ----
# Changing the (wrapped) text above causes the rows that never had an
# explicit/custom row height to auto-resize. Re-run AutoFit on those rows
# so they fall back to the sheet's default height, exactly like before the
# edit (rows that already had a custom height are left alone by AutoFit).
$sheets[1].Rows.Item(10).AutoFit()
$sheets[2].Rows.Item(10).AutoFit()
$sheets[3].Rows.Item(10).AutoFit()
$sheets[3].Rows.Item(11).AutoFit()
$sheets[4].Rows.Item(10).AutoFit()

# ---- Row height updates (matching the revised documentation layout) ----
$sheets[2].Rows.Item(11).RowHeight = 38.25
$sheets[2].Rows.Item(12).RowHeight = 51.0
$sheets[3].Rows.Item(12).RowHeight = 126.75
$sheets[4].Rows.Item(11).RowHeight = 28.5

# ---- Vertical alignment fix: UC 02.3 sheet, cell C12 (bottom -> center) ----
$sheets[3].Range("C12").VerticalAlignment = -4108  # xlCenter

# ---- Page margins for sheet1 (UC 02.1 Cadastrar Cliente) ----
# Excel COM margins are expressed in points; 1 inch = 72 points.
$sheets[1].PageSetup.TopMargin = 0.75 * 72
$sheets[1].PageSetup.BottomMargin = 0.75 * 72
$sheets[1].PageSetup.LeftMargin = 0.25 * 72
$sheets[1].PageSetup.RightMargin = 0.25 * 72

Write-Host "Edits applied successfully."
